# day-12.xlsx: add hLookup formulas + rename Sheet1 -> hLookup
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the active sheet (was "Sheet1") to "hLookup" ---
$ws.Name = "hLookup"

# --- Row 7: HLOOKUP into the 2nd row of the lookup table ($D$15:$J$17) ---
$ws.Range("E7").Formula = "=HLOOKUP(E6,`$D`$15:`$J`$17,2,FALSE)"
$ws.Range("F7:I7").Formula = "=HLOOKUP(F6,`$D`$15:`$J`$17,2,FALSE)"
$ws.Range("J7").Formula = "=HLOOKUP(J6,`$D`$15:`$J`$17,2,FALSE)"

# --- Row 10: HLOOKUP into the 3rd row of the lookup table ($D$15:$J$17) ---
$ws.Range("E10:F10").Formula = "=HLOOKUP(E6,`$D`$15:`$J`$17,3,FALSE)"
$ws.Range("G10").Formula = "=HLOOKUP(G6,`$D`$15:`$J`$17,3,FALSE)"
$ws.Range("H10:J10").Formula = "=HLOOKUP(H6,`$D`$15:`$J`$17,3,FALSE)"

# --- New helper column C alongside the lookup table (rows 15:17) ---
$ws.Range("C15").Value = 1
$ws.Range("C16").Value = 2
$ws.Range("C17").Value = 3

# --- View state: zoom + selection on the hLookup sheet ---
$excel.ActiveWindow.Zoom = 176
$ws.Range("E10:G10").Select()
